$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Q (2020) -------------------------------------------------
# Row 3 header: year 2020, same formatting as P3 (header cell).
$ws.Range("P3").Copy() | Out-Null
$ws.Range("Q3").PasteSpecial(-4122) | Out-Null
$ws.Range("Q3").Value = 2020

# Rows 4-12: data cells, all "-" (no data available), matching the
# formatting used across the rest of the table for this row band (P5 style).
$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q4:Q12").PasteSpecial(-4122) | Out-Null
$ws.Range("Q4:Q12").Value = "-"

# Row 13: bottom (thick-bottom-border) row, matching P13 style.
$ws.Range("P13").Copy() | Out-Null
$ws.Range("Q13").PasteSpecial(-4122) | Out-Null
$ws.Range("Q13").Value = "-"

$excel.CutCopyMode = 0

# --- Selection, as recorded by the author's last save --------------------
$ws.Range("P17").Select() | Out-Null
